# Rename the header in B1 from "cases" to "cases1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "cases1"

# Move the active selection from B1 to B2
$ws.Range("B2").Select()
